# Commit: "Fruta / hortaliza, semanal"
# Insert a new weekly data row for Berenjena (Macroferia Regional de Talca) at
# row 45, pushing the previously-existing rows 45..145 down to 46..146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45 (shifts rows 45-145 -> 46-146)
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly record
$ws.Range("A45").Value2 = 5
$ws.Range("B45").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C45").Value2 = 'Maule'
$ws.Range("D45").Value2 = 44987
$ws.Range("E45").Value2 = 7
$ws.Range("F45").Value2 = 100112001
$ws.Range("G45").Value2 = 'Berenjena'
$ws.Range("H45").Value2 = 'Sin especificar'
$ws.Range("I45").Value2 = 'Primera'
$ws.Range("J45").Value2 = 200
$ws.Range("K45").Value2 = 8000
$ws.Range("L45").Value2 = 8000
$ws.Range("M45").Value2 = 8000
$ws.Range("N45").Value2 = '$/caja 50 unidades'
$ws.Range("O45").Value2 = 'Región del Maule'
$ws.Range("P45").Value2 = 160
$ws.Range("Q45").Value2 = 50
$ws.Range("R45").Value2 = 'Hortaliza'
